# Update "想去人数" (interested-count) figures captured in a later scrape.
# Sheet order (per workbook.xml): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Range("F3").Value  = 35
$ws1.Range("F4").Value  = 356
$ws1.Range("F5").Value  = 41
$ws1.Range("F6").Value  = 6137
$ws1.Range("F7").Value  = 679
$ws1.Range("F9").Value  = 34
$ws1.Range("F14").Value = 1060
$ws1.Range("F15").Value = 63
$ws1.Range("F17").Value = 307
$ws1.Range("F18").Value = 1384
$ws1.Range("F20").Value = 1036
$ws1.Range("F22").Value = 2103
$ws1.Range("F23").Value = 212
$ws1.Range("F24").Value = 54
$ws1.Range("F25").Value = 373
$ws1.Range("F27").Value = 3408

$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws2.Range("F4").Value  = 353
$ws2.Range("F5").Value  = 111
$ws2.Range("F8").Value  = 24
$ws2.Range("F9").Value  = 671
$ws2.Range("F14").Value = 88
$ws2.Range("F18").Value = 5
$ws2.Range("F19").Value = 366
$ws2.Range("F21").Value = 4068
$ws2.Range("F24").Value = 30
$ws2.Range("F25").Value = 166
$ws2.Range("F27").Value = 80
$ws2.Range("F29").Value = 200

$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws3.Range("F8").Value  = 1528
$ws3.Range("F12").Value = 698

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Range("F6").Value  = 1528
$ws4.Range("F9").Value  = 698
$ws4.Range("F11").Value = 35
$ws4.Range("F12").Value = 356
$ws4.Range("F13").Value = 41
$ws4.Range("F14").Value = 6137
$ws4.Range("F15").Value = 24
$ws4.Range("F16").Value = 679
$ws4.Range("F18").Value = 34
$ws4.Range("F25").Value = 88
$ws4.Range("F28").Value = 1060
$ws4.Range("F29").Value = 63
$ws4.Range("F31").Value = 307
$ws4.Range("F32").Value = 5
$ws4.Range("F33").Value = 366
$ws4.Range("F34").Value = 1384
$ws4.Range("F36").Value = 30
$ws4.Range("F37").Value = 166
$ws4.Range("F39").Value = 1036
$ws4.Range("F42").Value = 2103
$ws4.Range("F44").Value = 212
$ws4.Range("F45").Value = 54
$ws4.Range("F46").Value = 373
$ws4.Range("F48").Value = 3408
